$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2156.037
$ws.Range("J112").Value = 2260.52
$ws.Range("L112").Value = 6781.559999999999
$ws.Range("N112").Value = -8997.559999999999

$ws.Range("H129").Value = 1052.2174
$ws.Range("J129").Value = 1116.1
$ws.Range("L129").Value = 3348.3
$ws.Range("N129").Value = -13348.3

$ws.Range("H131").Value = 5659.231
$ws.Range("I131").Value = 2000
$ws.Range("J131").Value = 6077.4287
$ws.Range("K131").Value = 6000
$ws.Range("L131").Value = 18232.2861
$ws.Range("M131").Value = -960
$ws.Range("N131").Value = -28312.2861

$ws.Range("H132").Value = 2952.8235
$ws.Range("I132").Value = 2625.682
$ws.Range("J132").Value = 3552.5833
$ws.Range("K132").Value = 7877.045999999999
$ws.Range("L132").Value = 10657.7499
$ws.Range("M132").Value = -5347.045999999999
$ws.Range("N132").Value = -15717.7499

$ws.Range("H138").Value = 3452586
$ws.Range("I138").Value = 8001587.5
$ws.Range("J138").Value = 6372.394
$ws.Range("K138").Value = 24004762.5
$ws.Range("L138").Value = 19117.182
$ws.Range("M138").Value = -23999622.5
$ws.Range("N138").Value = -29397.182

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H35").Value = 26118.5
$ws.Range("I35").Value = 26118.5
$ws.Range("K35").Value = 26118.5
$ws.Range("M35").Value = -25712.5

$ws.Range("H74").Value = 1052.4546
$ws.Range("I74").Value = 829.5
$ws.Range("J74").Value = 1320
$ws.Range("K74").Value = 829.5
$ws.Range("L74").Value = 1320
$ws.Range("M74").Value = 44.5
$ws.Range("N74").Value = -3068

$ws.Range("H77").Value = 1052.4546
$ws.Range("I77").Value = 829.5
$ws.Range("J77").Value = 1320
$ws.Range("K77").Value = 4147.5
$ws.Range("L77").Value = 6600
$ws.Range("M77").Value = 220.5
$ws.Range("N77").Value = -15336

$ws.Range("H122").Value = 1680.2593
$ws.Range("I122").Value = 1740.5238
$ws.Range("J122").Value = 1469.3334
$ws.Range("K122").Value = 5221.5714
$ws.Range("L122").Value = 4408.0002
$ws.Range("M122").Value = -2771.5714
$ws.Range("N122").Value = -9308.0002

$ws.Range("H123").Value = 34428.5
$ws.Range("J123").Value = 34428.5
$ws.Range("L123").Value = 34428.5
$ws.Range("N123").Value = -44228.5

$ws.Range("H131").Value = 26333
$ws.Range("J131").Value = 26333
$ws.Range("L131").Value = 26333
$ws.Range("N131").Value = -36413

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 7507.143
$ws.Range("I22").Value = 8691.666999999999
$ws.Range("J22").Value = 400
$ws.Range("K22").Value = 8691.666999999999
$ws.Range("L22").Value = 400
$ws.Range("M22").Value = -8518.666999999999
$ws.Range("N22").Value = -746

$ws.Range("H53").Value = 23745
$ws.Range("I53").Value = 24490
$ws.Range("J53").Value = 23000
$ws.Range("K53").Value = 24490
$ws.Range("L53").Value = 23000
$ws.Range("M53").Value = -23916
$ws.Range("N53").Value = -24148

$ws.Range("H99").Value = 2160.0667
$ws.Range("I99").Value = 1233.6666
$ws.Range("J99").Value = 2391.6667
$ws.Range("K99").Value = 1233.6666
$ws.Range("L99").Value = 2391.6667
$ws.Range("M99").Value = 264.3334
$ws.Range("N99").Value = -5387.6667

$ws.Range("H107").Value = 37280.734
$ws.Range("I107").Value = 49546.453
$ws.Range("K107").Value = 49546.453
$ws.Range("M107").Value = -47626.453

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 129.66667
$ws.Range("I7").Value = 94.5
$ws.Range("K7").Value = 94.5
$ws.Range("M7").Value = 18.5

$ws.Range("H31").Value = 32261670
$ws.Range("I31").Value = 50003012
$ws.Range("K31").Value = 50003012
$ws.Range("M31").Value = -50002717

$ws.Range("H34").Value = 32261670
$ws.Range("I34").Value = 50003012
$ws.Range("K34").Value = 50003012
$ws.Range("M34").Value = -50002810

$ws.Range("H99").Value = 2104.0833
$ws.Range("I99").Value = 2247.2
$ws.Range("J99").Value = 1388.5
$ws.Range("K99").Value = 2247.2
$ws.Range("L99").Value = 1388.5
$ws.Range("M99").Value = -749.1999999999998
$ws.Range("N99").Value = -4384.5

$ws.Range("H126").Value = 2104.0833
$ws.Range("I126").Value = 2247.2
$ws.Range("J126").Value = 1388.5
$ws.Range("K126").Value = 6741.599999999999
$ws.Range("L126").Value = 4165.5
$ws.Range("M126").Value = -4271.599999999999
$ws.Range("N126").Value = -9105.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 67.14286
$ws.Range("J12").Value = 46.272728
$ws.Range("L12").Value = 138.818184
$ws.Range("N12").Value = -484.818184

$ws.Range("H69").Value = 991.93335
$ws.Range("I69").Value = 450
$ws.Range("K69").Value = 1350
$ws.Range("M69").Value = -539

$ws.Range("H72").Value = 991.93335
$ws.Range("I72").Value = 450
$ws.Range("K72").Value = 4050
$ws.Range("M72").Value = 6

$ws.Range("H113").Value = 684.2381
$ws.Range("I113").Value = 515.7143
$ws.Range("J113").Value = 768.5
$ws.Range("K113").Value = 1547.1429
$ws.Range("L113").Value = 2305.5
$ws.Range("M113").Value = 622.8571000000002
$ws.Range("N113").Value = -6645.5

$ws.Range("H131").Value = 866.21
$ws.Range("I131").Value = 466.66666
$ws.Range("J131").Value = 891.71277
$ws.Range("K131").Value = 1399.99998
$ws.Range("L131").Value = 2675.13831
$ws.Range("M131").Value = 3640.00002
$ws.Range("N131").Value = -12755.13831

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws.Range("H123").Value = 10393
$ws.Range("J123").Value = 10393
$ws.Range("L123").Value = 10393
$ws.Range("N123").Value = -15293

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 200
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 200
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -790

$ws.Range("H27").Value = 200
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 200
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 200
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -414

$ws.Range("H40").Value = 11730
$ws.Range("I40").Value = 16095
$ws.Range("K40").Value = 16095
$ws.Range("M40").Value = -15959

$ws.Range("H122").Value = 13894022
$ws.Range("I122").Value = 20838992
$ws.Range("J122").Value = 4083.3333
$ws.Range("K122").Value = 62516976
$ws.Range("L122").Value = 12249.9999
$ws.Range("M122").Value = -62514526
$ws.Range("N122").Value = -17149.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3500
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 3500
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 3500
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -4748

$ws.Range("H65").Value = 3500
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 3500
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 17500
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -23740

$ws.Range("H123").Value = 40442
$ws.Range("J123").Value = 40442
$ws.Range("L123").Value = 40442
$ws.Range("N123").Value = -50242

$ws.Range("H132").Value = 1895.4783
$ws.Range("I132").Value = 2217.5557
$ws.Range("J132").Value = 736
$ws.Range("K132").Value = 6652.6671
$ws.Range("L132").Value = 2208
$ws.Range("M132").Value = -4122.6671
$ws.Range("N132").Value = -7268
